# Updated Argent prices in Excel
# Appends a new row (45) to each price sheet with date 2025-04-15 and the
# latest price carried forward from row 44.

$wb = $excel.ActiveWorkbook

$sheetData = @(
    @{ Name = "N-Dense";                   Price = "40" },
    @{ Name = "N-Type";                    Price = "41" },
    @{ Name = "N-type Wafer";              Price = "1.25" },
    @{ Name = "Cell Topcon 183mm";         Price = "0.3" },
    @{ Name = "Module Topcon 183mm";       Price = "0.1" },
    @{ Name = "Silver Rear_side";          Price = "5,293" },
    @{ Name = "Silver Busbar front-side";  Price = "7,924" },
    @{ Name = "Silver finger front-side";  Price = "7,974" },
    @{ Name = "USD_CNY";                   Price = "7.3208" }
)

foreach ($entry in $sheetData) {
    $ws = $wb.Worksheets.Item($entry.Name)

    # Force the new cells to be stored as text (matching every other row in
    # these sheets, which are all inline/text strings, not real dates or
    # numbers), then strip the number format back off so no extra cell
    # style is introduced.
    $newRange = $ws.Range("A45:B45")
    $newRange.NumberFormat = "@"
    $ws.Cells.Item(45, 1).Value = "2025-04-15"
    $ws.Cells.Item(45, 2).Value = $entry.Price
    $newRange.ClearFormats()
}
